# Fix "1 John x:y" references in column A that were accidentally
# concatenated with a stray "16" (verse-count leftover) suffix, e.g.
# "1 John 1:416" -> "1 John 1:4". Only rows 2..79 (below the header)
# are affected; column B ("Text") is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.EndsWith("16")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 2)
    }
}
